$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$nl = [char]11

$c = $t.Cell(1, 1)
$c.Range.Text = "58 x 22" + $nl + "  2    2" + $nl + "  ----" + $nl + "5|    |" + $nl + "8|    |"

$c = $t.Cell(1, 2)
$c.Range.Text = "89 x 68" + $nl + "  6    8" + $nl + "  ----" + $nl + "8|    |" + $nl + "9|    |"

$c = $t.Cell(1, 3)
$c.Range.Text = "39 x 36" + $nl + "  3    6" + $nl + "  ----" + $nl + "3|    |" + $nl + "9|    |"

$c = $t.Cell(2, 1)
$c.Range.Text = "14 x 33" + $nl + "  3    3" + $nl + "  ----" + $nl + "1|    |" + $nl + "4|    |"

$c = $t.Cell(2, 2)
$c.Range.Text = "15 x 35" + $nl + "  3    5" + $nl + "  ----" + $nl + "1|    |" + $nl + "5|    |"

$c = $t.Cell(2, 3)
$c.Range.Text = "26 x 16" + $nl + "  1    6" + $nl + "  ----" + $nl + "2|    |" + $nl + "6|    |"

$c = $t.Cell(3, 1)
$c.Range.Text = "13 x 10" + $nl + "  1    0" + $nl + "  ----" + $nl + "1|    |" + $nl + "3|    |"

$c = $t.Cell(3, 2)
$c.Range.Text = "27 x 74" + $nl + "  7    4" + $nl + "  ----" + $nl + "2|    |" + $nl + "7|    |"

$c = $t.Cell(3, 3)
$c.Range.Text = "31 x 31" + $nl + "  3    1" + $nl + "  ----" + $nl + "3|    |" + $nl + "1|    |"

$c = $t.Cell(4, 1)
$c.Range.Text = "46 x 45" + $nl + "  4    5" + $nl + "  ----" + $nl + "4|    |" + $nl + "6|    |"

$c = $t.Cell(4, 2)
$c.Range.Text = "35 x 64" + $nl + "  6    4" + $nl + "  ----" + $nl + "3|    |" + $nl + "5|    |"

$c = $t.Cell(4, 3)
$c.Range.Text = "18 x 42" + $nl + "  4    2" + $nl + "  ----" + $nl + "1|    |" + $nl + "8|    |"

$c = $t.Cell(5, 1)
$c.Range.Text = "91 x 88" + $nl + "  8    8" + $nl + "  ----" + $nl + "9|    |" + $nl + "1|    |"

$c = $t.Cell(5, 2)
$c.Range.Text = "45 x 39" + $nl + "  3    9" + $nl + "  ----" + $nl + "4|    |" + $nl + "5|    |"

$c = $t.Cell(5, 3)
$c.Range.Text = "91 x 15" + $nl + "  1    5" + $nl + "  ----" + $nl + "9|    |" + $nl + "1|    |"
